$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ip_address_list": swap rows 10 and 11 (527_Teijin moves up to row 10,
# 474 B_Austin moves down to row 11) and clear the "favorite" flag (col E) on
# rows 11 and 12.
# Use Range.Copy (cell-to-cell copy) instead of direct Value assignment so
# that writing the multi-line notes in column D does not trigger the
# "autofit row height on text entry" behaviour (which would add ht/
# customHeight attributes that are not present in the target file).
# A distant scratch row (far below any real data) is used as temporary
# storage and removed afterwards with Rows.Delete so no trace is left behind.
# ---------------------------------------------------------------------------
$wsIp = $wb.Worksheets.Item("ip_address_list")

$wsIp.Range("A10:D10").Copy($wsIp.Range("A200:D200"))   # stash 474 B_Austin (row10) in scratch row
$wsIp.Range("A11:D11").Copy($wsIp.Range("A10:D10"))      # 527_Teijin -> row10
$wsIp.Range("A200:D200").Copy($wsIp.Range("A11:D11"))    # 474 B_Austin -> row11
$wsIp.Rows.Item(200).Delete()                            # remove scratch row entirely

$wsIp.Range("E11").Value2 = 0
$wsIp.Range("E12").Value2 = 0

# ---------------------------------------------------------------------------
# Sheet "ip_address_fav_list": remove the "474 B_Austin" row (old row 3) and
# the "515_ZF Stara Boleslav" row (old row 5). After deleting old row 3,
# "527_Teijin" (old row 4) becomes row 3, and "515_ZF..." becomes row 4, which
# is then also removed, leaving just 3 rows total.
# ---------------------------------------------------------------------------
$wsFav = $wb.Worksheets.Item("ip_address_fav_list")

$wsFav.Rows.Item(3).Delete()
$wsFav.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Sheet "Settings": update default interface setting value.
# ---------------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B1").Value2 = 6

# ---------------------------------------------------------------------------
# Sheet "projects_bin2" (hidden backup/bin sheet): populate rows 3 and 4 with
# data moved out of the other sheets (515_ZF Stara Boleslav row, and the
# 518_Valeo disk-list entry).
# ---------------------------------------------------------------------------
$wsBin = $wb.Worksheets.Item("projects_bin2")
$wsDisk = $wb.Worksheets.Item("disk_list")

# Row 3: same data as ip_address_list row 12 / old ip_address_fav_list row 5
$wsIp.Range("A12:D12").Copy($wsBin.Range("A3:D3"))
$wsBin.Range("E3").Value2 = 1

# Row 4: same data as disk_list row 5 (518_Valeo). Columns A-E copied as-is.
$wsDisk.Range("A5:E5").Copy($wsBin.Range("A4:E4"))

# Column F needs the same text as disk_list!F5, but with a stray control
# character removed. Fix it up in a distant scratch cell first (so the
# in-place text edit's row-height side effect lands on a throwaway row),
# then copy the corrected value into place and discard the scratch row.
$wsDisk.Range("F5").Copy($wsDisk.Range("A500"))
$wsDisk.Range("A500").Characters(9,1).Text = ""
$wsDisk.Range("A500").Copy($wsBin.Range("F4"))
$wsDisk.Rows.Item(500).Delete()
